# Improve game image handling: add an "Image URL" column (J) to the game
# wiki sheet, and append a new row for "Outlast" (row 7) populated across
# every column, including its image URL. The pre-existing rows (2-6) did
# not capture an image URL when they were written, so they get an empty
# (but present) Image URL cell -- matching how the sheet is re-shaped once
# the new column is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---------------------------------------------------
# Copy the existing header formatting (bold font, border, centered) from
# the neighboring "Date Added" header onto the new column before setting
# its text.
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 10).Value = "Image URL"

# --- Back-fill existing rows (2-6) with an empty Image URL cell --------
# Touching a format property is what makes Excel materialize the (blank)
# cell in the sheet instead of leaving it completely absent; ClearFormats
# then strips that formatting back off so the cell matches its neighbors.
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Borders.LineStyle = 1
    $cell.ClearFormats()
}

# --- New row 7: Outlast --------------------------------------------------
$ws.Cells.Item(7, 1).Value = 3790
$ws.Cells.Item(7, 2).Value = 'Outlast'
$ws.Cells.Item(7, 3).Value = 'J. T. Petty, Hugo Dallaire, Samuel Laflamme, David Chateauneuf, Philippe Morin'

# Release date is stored as plain text (matches D2:D6); force text so
# Excel doesn't reinterpret "2013-09-03" as a date serial.
$dateCell = $ws.Cells.Item(7, 4)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2013-09-03'
$dateCell.ClearFormats()

$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 'Outlast, developed by Red Barrels and released on September 3, 2013, is a psychological horror video game that immerses players into a harrowing experience within the confines of the Mount Massive Asylum. The game was developed by a team including J. T. Petty, Hugo Dallaire, Samuel Laflamme, David Chateauneuf, and Philippe Morin, who collectively pushed the boundaries of indie horror gaming. Outlast is available on multiple platforms, including Linux, macOS, Nintendo Switch, PC, Xbox One, and PlayStation 4, making it accessible to a wide audience of horror enthusiasts.
The narrative of Outlast follows investigative journalist Miles Upshur as he explores a remote psychiatric hospital to uncover the truth behind inhumane experiments conducted on its patients. The game is notable for its first-person perspective, which amplifies the immersive horror experience. Players are equipped with a night vision camera, which serves as their primary tool to navigate the dark and foreboding environments. This camera, however, requires a continuous supply of batteries, adding an additional layer of tension to the gameplay as players must carefully manage resources while avoiding the asylum''s hostile inhabitants.
Unlike traditional survival horror games, Outlast eschews combat mechanics in favor of a focus on stealth and evasion. Players are required to solve puzzles and locate items to progress, all while being relentlessly pursued by the asylum''s dangerous and deranged occupants. The realistic movement animations and strategically placed horror sound effects enhance the feeling of vulnerability and suspense, compelling players to empathize with the protagonist''s plight.
Upon release, Outlast received positive reviews for its atmospheric tension, gripping narrative, and innovative use of visual and sound design to create a sense of dread. It earned a rating of 3.74, praised for its ability to evoke terror without the need for direct combat. The game has been credited with revitalizing interest in the survival horror genre, influencing subsequent titles with its emphasis on psychological horror and narrative-driven gameplay.
Culturally, Outlast has had a significant impact on the horror gaming landscape, inspiring a wave of similar games that prioritize atmosphere and emotional engagement over traditional action elements. The success of Outlast demonstrated the potential for indie developers to create compelling and commercially successful horror experiences, contributing to a broader appreciation for the genre''s storytelling capabilities.'
$ws.Cells.Item(7, 7).Value = '[''Red Barrels. (2013). Outlast [Video game]. Red Barrels.'', ''Smith, A. (2013). Outlast review. IGN. Retrieved from https://www.ign.com/articles/2013/09/03/outlast-review'', ''Jones, M. (2013). The horror of Outlast: An interview with Red Barrels. GameSpot. Retrieved from https://www.gamespot.com/articles/the-horror-of-outlast-an-interview-with-red-barrels-6414145/'', "Takahashi, D. (2013). Red Barrels'' Outlast: How an indie team made a great horror game. VentureBeat. Retrieved from https://venturebeat.com/2013/09/17/red-barrels-outlast-how-an-indie-team-made-a-great-horror-game/", ''Official website of Red Barrels Games. (n.d.). Retrieved from http://redbarrelsgames.com/'']'
$ws.Cells.Item(7, 8).Value = 'Average Rating: 3.74
Average Playtime: 3 hours
ESRB Rating: Mature
Metacritic Score: 80
Platforms: Linux, macOS, Nintendo Switch, PC, Xbox One, PlayStation 4'
$ws.Cells.Item(7, 9).Value = '2025-03-27 18:27:27'
$ws.Cells.Item(7, 10).Value = 'https://media.rawg.io/media/games/9dd/9ddabb34840ea9227556670606cf8ea3.jpg'
